$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'Dr. Rana Abo-Zaid, Dr. Nourhan Mahmoud, Dr. Shimaa Ahmad Mekki, Dr. Majorelle Magdy, Dr. Servinaz Sayed Mohammad'
$ws.Range("G4").Value = 'Dr. Alshimaa Atef, Dr. Shimaa Ahmad Mekki, Dr. Menna tuâ€™Allah Medhat, Dr. Heba Mahmoud Ali, Dr. Amira Sobhy, Dr. Hend Mahmoud'
$ws.Range("G5").Value = 'Dr. Nada Gouda, Dr. Fatma Elhady, Dr. Menna tu''Alllah Mohammad, Dr. Abeer Ragab'
$ws.Range("G6").Value = 'Dr. Nada Mohammad, Dr. Kerelos Zareef'
$ws.Range("G8").Value = 'Dr. Aya Saeed, Dr. Amal Awwad'
$ws.Range("G10").Value = 'Dr. Marina Youhanna, Dr. Arwa Al-Sayed, Dr. Amany Raafat, Dr. Esraa Mostafa, Dr. Madeha Saeed, Dr. Maryam Ahmad'
$ws.Range("G11").Value = 'Dr. Alaa Ashraf, Dr. Sarah Mahdy'
$ws.Range("G12").Value = 'Dr. Sarah Mahdy, Dr. Nouran Mahmoud'
$ws.Range("G18").Value = 'Dr. Shorok Mohammad, Dr. Yasmin, Dr. Remon, Dr. Aya Hanafy'
$ws.Range("G19").Value = 'Dr. Wafaa Ebida, Dr. Remon, Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Nardine, Dr. Naema Gomaa, Dr. Monica, Dr. Yassmen Ahmad, Dr. Maryam Ashraf'
$ws.Range("G20").Value = 'Dr. Wafaa Ebida, Dr. Remon, Dr. Aya Hanafy, Dr. Nardine, Dr. Marina Sorial, Dr. Yassmen Ahmad, Dr. Youstina Magdy'
$ws.Range("G21").Value = 'Dr. Shorok Mohammad, Dr. Yassmen Ahmad, Dr. Yasmin, Dr. Monica'
$ws.Range("G22").Value = 'Dr. Wafaa Ebida, Dr. Naema Gomaa, Dr. Remon, Dr. Monica'
$ws.Range("G23").Value = 'Dr. Wafaa Ebida, Dr. Yassmen Ahmad'
$ws.Range("G24").Value = 'Dr. Rana Abo-Zaid, Dr. Nourhan Mahmoud, Dr. Shimaa Ahmad Mekki, Dr. Majorelle Magdy, Dr. Servinaz Sayed Mohammad'
$ws.Range("G25").Value = 'Dr. Manar Montaser, Dr. Gehan Adel, Administrator, Dr. Alshimaa Atef'
$ws.Range("G26").Value = 'Dr. Alshimaa Atef, Dr. Shimaa Ahmad Mekki, Dr. Menna tuâ€™Allah Medhat, Dr. Heba Mahmoud Ali, Dr. Amira Sobhy, Dr. Hend Mahmoud'
$ws.Range("G27").Value = 'Dr. Nada Gouda, Dr. Fatma Elhady, Dr. Menna tu''Alllah Mohammad, Dr. Abeer Ragab'
$ws.Range("G28").Value = 'Dr. Nada Mohammad, Dr. Kerelos Zareef'
$ws.Range("G30").Value = 'Dr. Aya Saeed, Dr. Amal Awwad'
$ws.Range("G32").Value = 'Dr. Marina Youhanna, Dr. Arwa Al-Sayed, Dr. Amany Raafat, Dr. Esraa Mostafa, Dr. Madeha Saeed, Dr. Maryam Ahmad'
$ws.Range("G33").Value = 'Dr. Alaa Ashraf, Dr. Sarah Mahdy'
$ws.Range("G34").Value = 'Dr. Sarah Mahdy, Dr. Nouran Mahmoud'
$ws.Range("G40").Value = 'Dr. Shorok Mohammad, Dr. Yasmin, Dr. Remon, Dr. Aya Hanafy'
$ws.Range("G41").Value = 'Dr. Wafaa Ebida, Dr. Remon, Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Nardine, Dr. Naema Gomaa, Dr. Monica, Dr. Yassmen Ahmad, Dr. Maryam Ashraf'
$ws.Range("G42").Value = 'Dr. Wafaa Ebida, Dr. Remon, Dr. Aya Hanafy, Dr. Nardine, Dr. Marina Sorial, Dr. Yassmen Ahmad, Dr. Youstina Magdy'
$ws.Range("G43").Value = 'Dr. Shorok Mohammad, Dr. Yassmen Ahmad, Dr. Yasmin, Dr. Monica'
$ws.Range("G44").Value = 'Dr. Wafaa Ebida, Dr. Naema Gomaa, Dr. Remon, Dr. Monica'
$ws.Range("G45").Value = 'Dr. Wafaa Ebida, Dr. Yassmen Ahmad'
$ws.Range("G46").Value = 'Dr. Shimaa Ahmad Mekki, Dr. Hend Mahmoud, Dr. Nahla Nagiub, Dr. Nourhan Mahmoud'
$ws.Range("G48").Value = 'Dr. Nahla Nagiub, Dr. Nourhan Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad'
$ws.Range("G52").Value = 'Dr. Mariam Nour El-Din, Dr. Shimaa Ashraf'
$ws.Range("G54").Value = 'Dr. Marina Youhanna, Dr. Arwa Al-Sayed, Dr. Yasmeena Fattoh, Dr. Marwa Mustafa, Dr. Mai Mustafa, Dr. Eman M. Abo-Sakaya, Dr. Amany Raafat, Dr. Basma Hamed, Dr. Madeha Saeed, Dr. Merna Said, Dr. Maryam Ahmad'
$ws.Range("G58").Value = 'Dr. Afaf Abdallah, Dr. Amr Saeed'
$ws.Range("G60").Value = 'Dr. Nancy Abd Al-Shafy, Dr. Amr Saeed'
$ws.Range("G62").Value = 'Dr. Shorok Mohammad, Dr. Wafaa Ebida, Dr. Yassmen Ahmad, Dr. Aya Hanafy'
$ws.Range("G63").Value = 'Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Maryam Ashraf'
$ws.Range("G65").Value = 'Dr. Remon, Dr. Aya Hanafy, Dr. Salma Hassan, Dr. Nardine, Dr. Ola Abd Al-Fattah, Dr. Shorok Mohammad, Dr. Eman Samir Gabry'
$ws.Range("G66").Value = 'Dr. Aya Hanafy, Dr. Marina Sorial, Dr. Eman Mohammad Al, Dr. Maryam Ashraf, Dr. Monica'
$ws.Range("G68").Value = 'Dr. Shimaa Ahmad Mekki, Dr. Hend Mahmoud, Dr. Nahla Nagiub, Dr. Nourhan Mahmoud'
$ws.Range("G70").Value = 'Dr. Nahla Nagiub, Dr. Nourhan Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad'
$ws.Range("G74").Value = 'Dr. Mariam Nour El-Din, Dr. Shimaa Ashraf'
$ws.Range("G76").Value = 'Dr. Marina Youhanna, Dr. Arwa Al-Sayed, Dr. Yasmeena Fattoh, Dr. Marwa Mustafa, Dr. Mai Mustafa, Dr. Eman M. Abo-Sakaya, Dr. Amany Raafat, Dr. Basma Hamed, Dr. Madeha Saeed, Dr. Merna Said, Dr. Maryam Ahmad'
$ws.Range("G80").Value = 'Dr. Afaf Abdallah, Dr. Amr Saeed'
$ws.Range("G82").Value = 'Dr. Nancy Abd Al-Shafy, Dr. Amr Saeed'
$ws.Range("G84").Value = 'Dr. Shorok Mohammad, Dr. Wafaa Ebida, Dr. Yassmen Ahmad, Dr. Aya Hanafy'
$ws.Range("G85").Value = 'Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Maryam Ashraf'
$ws.Range("G87").Value = 'Dr. Remon, Dr. Aya Hanafy, Dr. Salma Hassan, Dr. Nardine, Dr. Ola Abd Al-Fattah, Dr. Shorok Mohammad, Dr. Eman Samir Gabry'
$ws.Range("G88").Value = 'Dr. Aya Hanafy, Dr. Marina Sorial, Dr. Eman Mohammad Al, Dr. Maryam Ashraf, Dr. Monica'
$ws.Range("G90").Value = 'Dr. Manar Montaser, Dr. Shimaa Ahmad Mekki, Dr. Mohammad El-Tanany'
$ws.Range("G92").Value = 'Dr. Nahla Nagiub, Dr. Nourhan Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad'
$ws.Range("G93").Value = 'Dr. Fatma Elhady, Dr. Menna tu''Alllah Mohammad, Dr. Abeer Ragab, Dr. Amera Ahmad Saad'
$ws.Range("G96").Value = 'Dr. Sara Nabil, Dr. Mariam Nour El-Din, Dr. Amal Awwad, Dr. Nourhan Mohammad'
$ws.Range("G98").Value = 'Dr. Marina Youhanna, Dr. Arwa Al-Sayed, Dr. Yasmeena Fattoh, Dr. Marwa Mustafa, Dr. Mai Mustafa, Dr. Eman M. Abo-Sakaya, Dr. Amany Raafat, Dr. Basma Hamed, Dr. Madeha Saeed, Dr. Merna Said, Dr. Maryam Ahmad'
$ws.Range("G104").Value = 'Dr. Nancy Abd Al-Shafy, Dr. Amr Saeed'
$ws.Range("G106").Value = 'Dr. Wafaa Ebida, Dr. Remon, Dr. Neveen Nashaat, Dr. Nardine, Dr. Youstina Magdy, Dr. Monica'
$ws.Range("G107").Value = 'Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Aya Hanafy, Dr. Yassmen Ahmad, Dr. Maryam Ashraf, Dr. Monica'
$ws.Range("G108").Value = 'Dr. Wafaa Ebida, Dr. Remon, Dr. Aya Hanafy, Dr. Nardine, Dr. Marina Sorial, Dr. Yassmen Ahmad, Dr. Youstina Magdy'
$ws.Range("G110").Value = 'Dr. Wafaa Ebida, Dr. Monica, Dr. Yassmen Ahmad'
$ws.Range("G111").Value = 'Dr. Naema Gomaa, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Marina Atef'
$ws.Range("G112").Value = 'Dr. Manar Montaser, Dr. Shimaa Ahmad Mekki, Dr. Mohammad El-Tanany'
$ws.Range("G114").Value = 'Dr. Nahla Nagiub, Dr. Nourhan Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad'
$ws.Range("G115").Value = 'Dr. Fatma Elhady, Dr. Menna tu''Alllah Mohammad, Dr. Abeer Ragab, Dr. Amera Ahmad Saad'
$ws.Range("G118").Value = 'Dr. Sara Nabil, Dr. Mariam Nour El-Din, Dr. Amal Awwad, Dr. Nourhan Mohammad'
$ws.Range("G120").Value = 'Dr. Marina Youhanna, Dr. Arwa Al-Sayed, Dr. Yasmeena Fattoh, Dr. Marwa Mustafa, Dr. Mai Mustafa, Dr. Eman M. Abo-Sakaya, Dr. Amany Raafat, Dr. Basma Hamed, Dr. Madeha Saeed, Dr. Merna Said, Dr. Maryam Ahmad'
$ws.Range("G126").Value = 'Dr. Nancy Abd Al-Shafy, Dr. Amr Saeed'
$ws.Range("G128").Value = 'Dr. Wafaa Ebida, Dr. Remon, Dr. Neveen Nashaat, Dr. Nardine, Dr. Youstina Magdy, Dr. Monica'
$ws.Range("G129").Value = 'Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Aya Hanafy, Dr. Yassmen Ahmad, Dr. Maryam Ashraf, Dr. Monica'
$ws.Range("G130").Value = 'Dr. Wafaa Ebida, Dr. Remon, Dr. Aya Hanafy, Dr. Nardine, Dr. Marina Sorial, Dr. Yassmen Ahmad, Dr. Youstina Magdy'
$ws.Range("G131").Value = 'Dr. Nardine, Dr. Marina Atef'
$ws.Range("G132").Value = 'Dr. Wafaa Ebida, Dr. Monica, Dr. Yassmen Ahmad'
$ws.Range("G133").Value = 'Dr. Naema Gomaa, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Marina Atef'
$ws.Range("G134").Value = 'Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Hend Mahmoud'
$ws.Range("G137").Value = 'Dr. Nada Gouda, Dr. Fatma Elhady, Dr. Menna tu''Alllah Mohammad, Dr. Abeer Ragab'
$ws.Range("G140").Value = 'Dr. Aya Saeed, Dr. Amal Awwad'
$ws.Range("G142").Value = 'Dr. Marwa Mustafa, Dr. Yasmeena Fattoh, Dr. Basma Hamed, Dr. Esraa Mostafa, Dr. Merna Said'
$ws.Range("G150").Value = 'Dr. Wafaa Ebida, Dr. Remon, Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Nardine, Dr. Naema Gomaa, Dr. Monica, Dr. Yassmen Ahmad, Dr. Maryam Ashraf'
$ws.Range("G151").Value = 'Dr. Wafaa Ebida, Dr. Marina Atef, Dr. Yassmen Ahmad, Dr. Monica'
$ws.Range("G152").Value = 'Dr. Wafaa Ebida, Dr. Marina Atef'
$ws.Range("G153").Value = 'Dr. Aya Hanafy, Dr. Marina Sorial, Dr. Eman Mohammad Al, Dr. Maryam Ashraf, Dr. Monica'
$ws.Range("G154").Value = 'Dr. Wafaa Ebida, Dr. Naema Gomaa, Dr. Remon, Dr. Salma Hassan'
$ws.Range("G155").Value = 'Dr. Naema Gomaa, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Marina Atef'
$ws.Range("G156").Value = 'Dr. Alshimaa Atef, Dr. Manar Montaser, Dr. Menna tuâ€™Allah Medhat, Dr. Majorelle Magdy, Dr. Mohammad El-Tanany'
$ws.Range("G159").Value = 'Dr. Nada Gouda, Dr. Fatma Elhady, Dr. Menna tu''Alllah Mohammad, Dr. Abeer Ragab'
$ws.Range("G162").Value = 'Dr. Aya Saeed, Dr. Amal Awwad'
$ws.Range("G164").Value = 'Dr. Marwa Mustafa, Dr. Yasmeena Fattoh, Dr. Basma Hamed, Dr. Esraa Mostafa, Dr. Merna Said'
$ws.Range("G165").Value = 'Dr. Sarah Mahdy, Dr. Nouran Mahmoud'
$ws.Range("G172").Value = 'Dr. Wafaa Ebida, Dr. Remon, Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Nardine, Dr. Naema Gomaa, Dr. Monica, Dr. Yassmen Ahmad, Dr. Maryam Ashraf'
$ws.Range("G173").Value = 'Dr. Wafaa Ebida, Dr. Marina Atef, Dr. Yassmen Ahmad, Dr. Monica'
$ws.Range("G174").Value = 'Dr. Wafaa Ebida, Dr. Marina Atef'
$ws.Range("G175").Value = 'Dr. Aya Hanafy, Dr. Marina Sorial, Dr. Eman Mohammad Al, Dr. Maryam Ashraf, Dr. Monica'
$ws.Range("G176").Value = 'Dr. Wafaa Ebida, Dr. Naema Gomaa, Dr. Remon, Dr. Salma Hassan'
$ws.Range("G177").Value = 'Dr. Naema Gomaa, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Marina Atef'
